# Apply updated cryptocurrency price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.529.27'
$ws.Range('E2').Value = '  +2.18%  '
$ws.Range('D3').Value = '''1.579.02'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('E4').Value = '  +0.62%  '
$ws.Range('D5').Value = '''212.43'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('D8').Value = '''46.79'
$ws.Range('E8').Value = '  +7.88%  '
$ws.Range('D9').Value = '''24.00'
$ws.Range('E9').Value = '  +4.08%  '
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '''1.803.39'
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').Value = '''1.590.14'
$ws.Range('E14').Value = '  +1.25%  '
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('E16').Value = '  -0.89%  '
$ws.Range('D17').Value = '''28.549.74'
$ws.Range('E17').Value = '  +2.31%  '
$ws.Range('D18').Value = '''62.45'
$ws.Range('D19').Value = '''229.49'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('E21').Value = '  -1.03%  '
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('D23').Value = '''3.94'
$ws.Range('E23').Value = '  -3.50%  '
$ws.Range('D24').Value = '''9.18'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D25').Value = '''2.03'
$ws.Range('E25').Value = '  +5.64%  '
$ws.Range('D26').Value = '''151.67'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('D28').Value = '''6.48'
$ws.Range('E28').Value = '  -1.16%  '
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('E32').Value = '  -1.33%  '
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range('E34').Value = '  +0.32%  '
$ws.Range('D35').Value = '''1.398.40'
$ws.Range('E35').Value = '  -0.66%  '
$ws.Range('E36').Value = '  -1.03%  '
$ws.Range('D37').Value = '''1.02'
$ws.Range('E37').Value = '  -3.00%  '
$ws.Range('E38').Value = '  +2.55%  '
$ws.Range('D39').Value = '''2.59'
$ws.Range('E39').Value = '  +7.09%  '
$ws.Range('E40').Value = '  -0.18%  '
$ws.Range('D41').Value = '''0.534'
$ws.Range('E41').Value = '  -0.66%  '
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('E43').Value = '  -1.12%  '
$ws.Range('D44').Value = '''5.62'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('D45').Value = '''1.85'
$ws.Range('E45').Value = '  +2.54%  '
$ws.Range('D46').Value = '''0.980'
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('D47').Value = '''62.69'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('D48').Value = '''1.714.71'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').Value = '''85.85'
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('E51').Value = '  -1.23%  '
